# Auto-generated edit script: updates Leve profit calculation sheets
# with refreshed market-price data across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# Sheet index 1, Row 100
$ws = $wb.Worksheets.Item(1)
$ws.Range("H100").Value = 3999.8572
$ws.Range("I100").Value = 3799.8
$ws.Range("J100").Value = 4500
$ws.Range("K100").Value = 3799.8
$ws.Range("L100").Value = 4500
$ws.Range("M100").Value = -3258.8
$ws.Range("N100").Value = -5582

# Sheet index 1, Row 137
$ws = $wb.Worksheets.Item(1)
$ws.Range("H137").Value = 64241.5
$ws.Range("I137").Value = 1558.125
$ws.Range("K137").Value = 4674.375
$ws.Range("M137").Value = -2124.375

# Sheet index 1, Row 138
$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value = 1823.5161
$ws.Range("I138").Value = 1199.7812
$ws.Range("J138").Value = 2488.8333
$ws.Range("K138").Value = 3599.3436
$ws.Range("L138").Value = 7466.499899999999
$ws.Range("M138").Value = 1540.6564
$ws.Range("N138").Value = -17746.4999

# Sheet index 1, Row 141
$ws = $wb.Worksheets.Item(1)
$ws.Range("H141").Value = 2659.9
$ws.Range("I141").Value = 2048.625
$ws.Range("K141").Value = 6145.875
$ws.Range("M141").Value = -965.875

# Sheet index 2, Row 32
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 19087.068
$ws.Range("I32").Value = 19413.158
$ws.Range("K32").Value = 19413.158
$ws.Range("M32").Value = -19126.158

# Sheet index 2, Row 61
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 2331.7805
$ws.Range("I61").Value = 1957.08
$ws.Range("J61").Value = 2917.25
$ws.Range("K61").Value = 1957.08
$ws.Range("L61").Value = 2917.25
$ws.Range("M61").Value = -1745.08
$ws.Range("N61").Value = -3341.25

# Sheet index 2, Row 74
$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 50002810
$ws.Range("I74").Value = 62503196
$ws.Range("K74").Value = 62503196
$ws.Range("M74").Value = -62502322

# Sheet index 2, Row 77
$ws = $wb.Worksheets.Item(2)
$ws.Range("H77").Value = 50002810
$ws.Range("I77").Value = 62503196
$ws.Range("K77").Value = 312515980
$ws.Range("M77").Value = -312511612

# Sheet index 2, Row 132
$ws = $wb.Worksheets.Item(2)
$ws.Range("H132").Value = 16213.371
$ws.Range("I132").Value = 1767.0435
$ws.Range("K132").Value = 5301.1305
$ws.Range("M132").Value = -2771.1305

# Sheet index 2, Row 136
$ws = $wb.Worksheets.Item(2)
$ws.Range("H136").Value = 2331.7805
$ws.Range("I136").Value = 1957.08
$ws.Range("J136").Value = 2917.25
$ws.Range("K136").Value = 5871.24
$ws.Range("L136").Value = 8751.75
$ws.Range("M136").Value = -3321.24
$ws.Range("N136").Value = -13851.75

# Sheet index 3, Row 26
$ws = $wb.Worksheets.Item(3)
$ws.Range("H26").Value = 22823.666
$ws.Range("I26").Value = 19235.5
$ws.Range("J26").Value = 30000
$ws.Range("K26").Value = 19235.5
$ws.Range("L26").Value = 30000
$ws.Range("M26").Value = -18943.5
$ws.Range("N26").Value = -30584

# Sheet index 3, Row 94
$ws = $wb.Worksheets.Item(3)
$ws.Range("H94").Value = 977.67346
$ws.Range("I94").Value = 951.9
$ws.Range("K94").Value = 951.9
$ws.Range("M94").Value = -500.9

# Sheet index 3, Row 99
$ws = $wb.Worksheets.Item(3)
$ws.Range("H99").Value = 1743.5
$ws.Range("I99").Value = 1176.25
$ws.Range("J99").Value = 2499.8333
$ws.Range("K99").Value = 1176.25
$ws.Range("L99").Value = 2499.8333
$ws.Range("M99").Value = 321.75
$ws.Range("N99").Value = -5495.8333

# Sheet index 3, Row 107
$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 639.625
$ws.Range("I107").Value = 602.8333
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 602.8333
$ws.Range("L107").Value = 750
$ws.Range("M107").Value = 1317.1667
$ws.Range("N107").Value = -4590

# Sheet index 3, Row 134
$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 38797.484
$ws.Range("I134").Value = 50642.137
$ws.Range("K134").Value = 151926.411
$ws.Range("M134").Value = -149391.411

# Sheet index 4, Row 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("H4").Value = 7000
$ws.Range("J4").Value = 7000
$ws.Range("L4").Value = 7000
$ws.Range("N4").Value = -7224

# Sheet index 4, Row 58
$ws = $wb.Worksheets.Item(4)
$ws.Range("H58").Value = 21104.08
$ws.Range("J58").Value = 72158.71000000001
$ws.Range("L58").Value = 72158.71000000001
$ws.Range("N58").Value = -72564.71000000001

# Sheet index 4, Row 132
$ws = $wb.Worksheets.Item(4)
$ws.Range("H132").Value = 18451.906
$ws.Range("I132").Value = 21209.576
$ws.Range("K132").Value = 63628.728
$ws.Range("M132").Value = -61098.728

# Sheet index 4, Row 134
$ws = $wb.Worksheets.Item(4)
$ws.Range("H134").Value = 1058.7693
$ws.Range("I134").Value = 989.45
$ws.Range("K134").Value = 2968.35
$ws.Range("M134").Value = -433.3500000000004

# Sheet index 4, Row 136
$ws = $wb.Worksheets.Item(4)
$ws.Range("H136").Value = 21104.08
$ws.Range("J136").Value = 72158.71000000001
$ws.Range("L136").Value = 216476.13
$ws.Range("N136").Value = -221576.13

# Sheet index 5, Row 129
$ws = $wb.Worksheets.Item(5)
$ws.Range("H129").Value = 334667.94
$ws.Range("I129").Value = 633.3333
$ws.Range("J129").Value = 418176.6
$ws.Range("K129").Value = 1899.9999
$ws.Range("L129").Value = 1254529.8
$ws.Range("M129").Value = 3100.0001
$ws.Range("N129").Value = -1264529.8

# Sheet index 5, Row 131
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 726.5
$ws.Range("J131").Value = 726.5
$ws.Range("L131").Value = 2179.5
$ws.Range("N131").Value = -12259.5

# Sheet index 5, Row 139
$ws = $wb.Worksheets.Item(5)
$ws.Range("H139").Value = 1615.7084
$ws.Range("I139").Value = 1103.6
$ws.Range("J139").Value = 2469.2222
$ws.Range("K139").Value = 3310.8
$ws.Range("L139").Value = 7407.6666
$ws.Range("M139").Value = 1829.2
$ws.Range("N139").Value = -17687.6666

# Sheet index 5, Row 140
$ws = $wb.Worksheets.Item(5)
$ws.Range("H140").Value = 1461.1904
$ws.Range("I140").Value = 1271.0526
$ws.Range("J140").Value = 3267.5
$ws.Range("K140").Value = 3813.1578
$ws.Range("L140").Value = 9802.5
$ws.Range("M140").Value = 1366.8422
$ws.Range("N140").Value = -20162.5

# Sheet index 5, Row 141
$ws = $wb.Worksheets.Item(5)
$ws.Range("H141").Value = 3732.0908
$ws.Range("I141").Value = 6765
$ws.Range("J141").Value = 1999
$ws.Range("K141").Value = 20295
$ws.Range("L141").Value = 5997
$ws.Range("M141").Value = -15115
$ws.Range("N141").Value = -16357

# Sheet index 6, Row 132
$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 41682.152
$ws.Range("I132").Value = 42145.46
$ws.Range("J132").Value = 40755.54
$ws.Range("K132").Value = 126436.38
$ws.Range("L132").Value = 122266.62
$ws.Range("M132").Value = -123906.38
$ws.Range("N132").Value = -127326.62

# Sheet index 7, Row 2
$ws = $wb.Worksheets.Item(7)
$ws.Range("H2").Value = 1178571.4
$ws.Range("I2").Value = 1230769.2
$ws.Range("J2").Value = 500000
$ws.Range("K2").Value = 1230769.2
$ws.Range("L2").Value = 500000
$ws.Range("M2").Value = -1230657.2
$ws.Range("N2").Value = -500224

# Sheet index 7, Row 132
$ws = $wb.Worksheets.Item(7)
$ws.Range("H132").Value = 1720.7742
$ws.Range("I132").Value = 1146.7222
$ws.Range("J132").Value = 2515.6155
$ws.Range("K132").Value = 3440.1666
$ws.Range("L132").Value = 7546.8465
$ws.Range("M132").Value = -910.1665999999996
$ws.Range("N132").Value = -12606.8465

# Sheet index 7, Row 136
$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value = 34085.734
$ws.Range("I136").Value = 34085.734
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 102257.202
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -99707.20199999999
$ws.Range("N136").ClearContents()

# Sheet index 8, Row 81
$ws = $wb.Worksheets.Item(8)
$ws.Range("H81").Value = 1349.4166
$ws.Range("I81").Value = 1611.375
$ws.Range("J81").Value = 825.5
$ws.Range("K81").Value = 3222.75
$ws.Range("L81").Value = 1651
$ws.Range("M81").Value = -2161.75
$ws.Range("N81").Value = -3773

# Sheet index 8, Row 84
$ws = $wb.Worksheets.Item(8)
$ws.Range("H84").Value = 1349.4166
$ws.Range("I84").Value = 1611.375
$ws.Range("J84").Value = 825.5
$ws.Range("K84").Value = 16113.75
$ws.Range("L84").Value = 8255
$ws.Range("M84").Value = -10809.75
$ws.Range("N84").Value = -18863

# Sheet index 8, Row 132
$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 1594.5454
$ws.Range("I132").Value = 838
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 2514
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = 16
$ws.Range("N132").Value = -20057

# Sheet index 8, Row 136
$ws = $wb.Worksheets.Item(8)
$ws.Range("H136").Value = 17858552
$ws.Range("I136").Value = 25642252
$ws.Range("J136").Value = 1826.1177
$ws.Range("K136").Value = 76926756
$ws.Range("L136").Value = 5478.3531
$ws.Range("M136").Value = -76924206
$ws.Range("N136").Value = -10578.3531

